$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 ("Changes in PWM can't be recognized/decoded") has been mitigated.
$ws.Range("F4").Value = "Mitigated"
$ws.Range("G4").Value = "M"
$ws.Range("H4").Value = "Differences have been made very noticable"

# Update the active selection to match the author's final cursor position.
$ws.Range("B4").Select()
